$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / period update
$ws.Range("E23").Value = "26.02.2026"
$ws.Range("C24").Value = "04.02.2026 bis 26.02.2026"

# Fix text label for rows 33 & 35
$ws.Range("C33").Value = "Unterricht ASA 8"
$ws.Range("C35").Value = "Unterricht ASA 8"

# Row 36
$ws.Range("A36").Value = "13.02.2026"
$ws.Range("B36").Value = 3.25
$ws.Range("C36").Value = "Unterricht ASA 9 Vertretung"

# Row 37
$ws.Range("A37").Value = "16.02.2026"
$ws.Range("B37").Value = 0.75
$ws.Range("C37").Value = "Vorbereitung für Unterricht: 0,75 Stunden"

# Row 38
$ws.Range("A38").Value = "17.02.2026"
$ws.Range("B38").Value = 3
$ws.Range("C38").Value = "Unterricht ASA 8"

# Row 39
$ws.Range("A39").Value = "18.02.2026"
$ws.Range("B39").Value = 0.75
$ws.Range("C39").Value = "Vorbereitung für Unterricht: 0,75 Stunden"

# Row 40
$ws.Range("A40").Value = "19.02.2026"
$ws.Range("B40").Value = 3
$ws.Range("C40").Value = "Unterricht ASA 8"

# Row 41
$ws.Range("A41").Value = "25.02.2026"
$ws.Range("B41").Value = 0.75
$ws.Range("C41").Value = "Vorbereitung für Unterricht: 0,75 Stunden"
$ws.Range("D41").Value = 19

# Row 42
$ws.Range("A42").Value = "26.02.2026"
$ws.Range("B42").Value = 3
$ws.Range("C42").Value = "Unterricht ASA 8"
$ws.Range("D42").Value = 19
